$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ------------------------------------------------------------------
# Section 1: replace the "Objects were used..." paragraph (class
# structure intro) with the new two-paragraph description of the
# City/Route classes.
# ------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute("Objects were used in order to make the project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$r1.Expand(4)
$xml1 = '<w:p><w:r><w:t>The data is broken down into two easily-defined classes: Cities and Routes.  These classes have been outfitted for ease-of-use.  They both have copy-constructors and assignment operator overloads to make assigning values much easier, as well as un-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>argumented</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>argumented</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> constructors for declaration.  The City class acts as a container for the ID and coordinates for a city.  It has a getter and setter for each data field.  The Route class contains a versatile array of Cities.  The array can be set using another array of Cities, the Route’s copy constructor, or using another Route with an assignment statement.  The main use of the Route class is its distance function, which returns the overall distance between the cities, including the distance between the last city and the first.</w:t></w:r></w:p><w:p><w:r><w:t>The decision to use object oriented programming was made very early on to simplify the organization of the routes and the cities contained in them.  The inclusion of the variety of methods to reassign a route makes using the classes very flexible, and speeds up debugging as well.  Since the classes are relatively basic, no significant problems were encountered during the creation or testing process.</w:t></w:r></w:p>'
$xml1 = $xml1.Replace("<w:p>", "<w:p " + $ns + ">")
[void]$r1.InsertXML($xml1)

# ------------------------------------------------------------------
# Sections 2-4: the merge-sort timing paragraph ("...was reduced to
# approximately 57 minutes") and the crossover-function paragraph
# (adds a lastRenderedPageBreak up front, and removes the mid-sentence
# one around "over cities ... with more connections").
# ------------------------------------------------------------------
$r2a = $d.Content
[void]$r2a.Find.Execute("For the sorting of our parents", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$r2a.Expand(4)

$r2b = $d.Content
[void]$r2b.Find.Execute("The crossover function that is used", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$r2b.Expand(4)

$r2 = $d.Range($r2a.Start, $r2b.End)
$xml2 = '<w:p><w:r><w:t xml:space="preserve">For the sorting of our parents, we use a merge sort due to its efficiency and </w:t></w:r><w:r><w:t xml:space="preserve">its </w:t></w:r><w:r><w:t xml:space="preserve">ability to </w:t></w:r><w:r><w:t>be parallelized</w:t></w:r><w:r><w:t>.  To start our proje</w:t></w:r><w:r><w:t>ct, we used a bubble sort</w:t></w:r><w:r><w:t>.  With 1000 parents the bubble sort version ran in approximately 13 hours.  With the implementation of a merge sort in serial, our run time</w:t></w:r><w:r><w:t xml:space="preserve"> was reduced to approximately 57</w:t></w:r><w:r><w:t xml:space="preserve"> minutes.  </w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">The crossover function that is used to generate the new solutions is based on the Edge Recombination Algorithm.  This algorithm works by merging two routes into a new route that mostly consists of connections that existed in either of the parents.  This algorithm also has an element of randomness in it, which can cause mutations through a few different methods.  The algorithm works by taking a union of both parent’s adjacency matrices.  This adjacency matrix represents any possibly </w:t></w:r><w:r><w:t>connection</w:t></w:r><w:r><w:t xml:space="preserve"> that can be followed from one city to another in either of the parent routes.  The route is then created by selecting the next city to travel to by valuing cities that have the least amount of connections </w:t></w:r><w:r><w:t xml:space="preserve">over cities with more connections.  If there is ever a point where there are no more possible connections within the parents that haven’t been used already for the current city then the next city is picked randomly from the remaining destinations that have not been travelled to – this is a mutation.  Mutations are pretty rare, but do occur.    </w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r></w:p>'
$xml2 = $xml2.Replace("<w:p>", "<w:p " + $ns + ">")
[void]$r2.InsertXML($xml2)

# ------------------------------------------------------------------
# Section 5: split the "merge sort will run in true serial..."
# paragraph so a lastRenderedPageBreak lands before "stunted serial
# run time".
# ------------------------------------------------------------------
$r3 = $d.Content
[void]$r3.Find.Execute("The merge sort will run in true serial", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$r3.Expand(4)
$xml3 = '<w:p><w:r><w:t xml:space="preserve">The merge sort will run in true serial if NUM_THREADS is placed to 1, however the parallel region around the Genetic Algorithm code might have some overhead in this situation which would result in a </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">stunted serial run time.  That is why it is suggested to turn off the OpenMP master switch inside of Visual Studio if a serial run is desired.  The instructions for turning off OpenMP in Visual Studio can be found in the User Manual. </w:t></w:r></w:p>'
$xml3 = $xml3.Replace("<w:p>", "<w:p " + $ns + ">")
[void]$r3.InsertXML($xml3)

# ------------------------------------------------------------------
# Section 6: merge the "Test runs showed..." paragraph's first two
# runs back together, dropping the lastRenderedPageBreak that used to
# sit before "the first generation".
# ------------------------------------------------------------------
$r4 = $d.Content
[void]$r4.Find.Execute("Test runs showed a small improvement", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$r4.Expand(4)
$xml4 = '<w:p><w:r><w:t xml:space="preserve">Test runs showed a small improvement in run-time in parallel versus serial.  The speed up averaged around 1.2-1.5 times faster when ran in parallel.  Because of the potential randomness in the creation of the first generation the best way to compare run times between solutions is to calculate iterations per second.  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>An iteration</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is defined as the process of creating new children and finding the best solution out of the new generation mixed with the old generation.  Since the number of iterations can vary from run to run, and this would influence total run time, it is better to compare how long the code took per iteration. </w:t></w:r></w:p>'
$xml4 = $xml4.Replace("<w:p>", "<w:p " + $ns + ">")
[void]$r4.InsertXML($xml4)
